$wb = $excel.ActiveWorkbook
$count = $wb.Worksheets.Count
$last = $wb.Worksheets.Item($count)
$new = $wb.Worksheets.Add($null, $last)
$new.Name = "magapoke_2025-12-10"

$headerRange = $new.Range("A1:B1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$new.Range("A1").Value = "rank"
$new.Range("B1").Value = 'title'

$titles = @(
    'ブルーロック'
    'WIND BREAKER'
    '信じていた仲間達にダンジョン奥地で殺されかけたがギフト『無限ガチャ』でレベル9999の仲間達を手に入れて元パーティーメンバーと世界に復讐＆『ざまぁ！』します！'
    'みいちゃんと山田さん'
    '東京卍リベンジャーズ'
    'ガチアクタ'
    'ベイビーステップ'
    'ギルティサークル'
    '島耕作'
    'FAIRY TAIL 100 YEARS QUEST'
    'イレギュラーズ'
    '薫る花は凛と咲く'
    '十字架のろくにん'
    '愛妻の裏アカ'
    'ハードワーカー中田'
    '南海トラフ巨大地震'
    '魔術ギルド総帥～生まれ変わって今更やり直す2度目の学院生活～'
    '黄昏町プリズナーズ'
    '転生したら第七王子だったので、気ままに魔術を極めます'
    '君が僕らを悪魔と呼んだ頃'
    '異世界ウォーキング'
    'GALAXIAS'
    '味方が弱すぎて補助魔法に徹していた宮廷魔法師、追放されて最強を目指す'
    'K-9~警視庁公安部公安第9課異能対策係~'
    '蒼く染めろ'
    '転生貴族、鑑定スキルで成り上がる～弱小領地を受け継いだので、優秀な人材を増やしていたら、最強領地になってた～'
    'アルキメデスの大戦'
    'ドラハチ'
    'さわらないで小手指くん'
    '魔女と傭兵'
    'おやすみ ふみさん'
    '異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～'
    'せいぶつ部の田辺くん'
    'グラぱらっ！'
    '屋根の下のアルテミス'
    'ひゃくえむ。'
    '黒猫と魔女の教室'
    '降り積もれ孤独な死よ'
    '幼馴染とはラブコメにならない'
    'ともだちづくり'
    'ナキナギ'
    '念願の悪役令嬢（ラスボス）の身体を手に入れたぞ！'
    'となりの黒川さん'
    'ハナバス　苔石花江のバスケ論'
    '追放された転生王子、『自動製作《オートクラフト》』スキルで領地を爆速で開拓し最強の村を作ってしまう〜最強クラフトスキルで始める、楽々領地開拓スローライフ〜'
    '食糧人類-Starving Anonymous-'
    '追放されなかった男　～二度目の人生は土下座から始まりました～'
    'いじめるヤバイ奴'
    'アオバノバスケ'
    '阿武ノーマル'
    '五輪の女神さま 〜なでしこ寮のメダルごはん〜'
    'イジらないで、長瀞さん'
    '最弱な僕は＜壁抜けバグ＞で成り上がる～壁をすり抜けたら、初回クリア報酬を無限回収できました！～'
    '皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～'
    'Destiny Unchain Online 〜吸血鬼少女となって、やがて『赤の魔王』と呼ばれるようになりました〜'
    '限界集落を脱村した錬金術士、都会で"最強"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～'
    'デッドアカウント'
    '春くらり'
    'MYS'
    'ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜'
    'Aランクパーティを離脱した俺は、元教え子たちと迷宮深部を目指す。'
    '東京卍リベンジャーズ～場地圭介からの手紙～'
    '普通の本はありません！'
    '劣等人の魔剣使い　スキルボードを駆使して最強に至る'
    '時々ボソッとロシア語でデレる隣のアーリャさん'
    'インフェクション'
    'お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！'
    '四十七大戦'
    '不遇職【鑑定士】が実は最強だった～奈落で鍛えた最強の【神眼】で無双する～'
    'リスナーに騙されてダンジョンの最下層から脱出RTAすることになった'
    'なれの果ての僕ら'
    'この世界がいずれ滅ぶことを、俺だけが知っている～モンスターが現れた世界で、死に戻りレベルアップ～'
    '我間乱 ―修羅―'
    'ダメスキル【自動機能】が覚醒しました～あれ、ギルドのスカウトの皆さん、俺を「いらない」って言ってませんでした？～'
    'ストーカー行為がバレて人生終了男'
    '君が監督！'
    'デスティニーラバーズ'
    'ハプスブルク家の華麗なる受難'
    '剣帝学院の魔眼賢者'
    '中華一番！極'
    '人間消失'
    '不遇職『鍛冶師』だけど最強です ～気づけば何でも作れるようになっていた男ののんびりスローライフ～'
    'ヒロインは絶望しました。'
    'シャングリラ・フロンティア～クソゲーハンター、神ゲーに挑まんとす～'
    'それがメイドのカンナです'
    '母という呪縛 娘という牢獄'
    '四刀流の最強配信者～やり込んだVRゲームの設定が現実世界に反映されたので、廃止予定だった戦闘職で無双します～'
    '恋ニ非ズ'
    '辺境の薬師、都でSランク冒険者となる～英雄村の少年がチート薬で無自覚無双〜'
    'ジュミドロ'
    '復讐の教科書'
    '可愛いだけじゃない式守さん'
    '彼女、お借りします'
    '東京ネオンスキャンダル'
    '触手魔術師の成り上がり'
    'お願い、脱がシて。'
    'お嬢様の僕'
    'はじめの一歩'
    'DAYS外伝'
    '魁の花巫女'
)

for ($i = 0; $i -lt $titles.Count; $i++) {
    $row = $i + 2
    $new.Cells.Item($row, 1).Value = $i + 1
    $new.Cells.Item($row, 2).Value = $titles[$i]
}